# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds the strikeout count for each outing. The values were
# regenerated from the underlying pitch-by-pitch log (using the actual "K"
# result code) instead of the old "Strike#" (total strikes thrown) figure
# used previously, so only column G changes row-by-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by worksheet row (column G is the 7th column)
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
